# Generate Report for Handback
# Updates the Handoff/Handback timestamps for the "4825cdc2-..." source file
# row in both the "zh-cn" and "de-de" language sheets, simulating a fresh
# report generation run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-24 06:49:00"
$zhcn.Range("K2").Value = "2016-08-24 06:49:27"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-24 06:49:10"
$dede.Range("K2").Value = "2016-08-24 06:49:35"
